# Revert "Ad hoc fix to pull_data":
#  - column A (date) goes back from inline-string "yyyy-mm-dd" text cells
#    to real numeric date serials formatted with a yyyy-mm-dd number format
#  - rows 52-62 regain their B (usphpi) / C (casusxam) figures that the
#    ad hoc fix had blanked out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial numbers for A2:A62 (1970-12-31 .. 2030-12-31)
$dates = @(25933,26298,26664,27029,27394,27759,28125,28490,28855,29220,29586,29951,30316,30681,31047,31412,31777,32142,32508,32873,33238,33603,33969,34334,34699,35064,35430,35795,36160,36525,36891,37256,37621,37986,38352,38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022,46387,46752,47118,47483,47848)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
}

# Apply the date number format to the (now numeric) date column, rows 2-62 only
# (row 1 keeps its existing bold/centered header style).
$ws.Range("A2:A62").NumberFormat = "yyyy-mm-dd"

# Restore the column width that Excel recorded once the column held real
# dates formatted with the custom number format.
$ws.Columns.Item(1).ColumnWidth = 19.8

# Restore the previously-blanked usphpi (B) / casusxam (C) values for 2020-2030
$restoredB = @(257.6916666666667,262.6638888888889,264.0951851851852,261.4835802469136,262.7475514403292,262.775438957476,262.3355235482396,262.6195046486816,262.5768223847991,262.5106168605735,262.5689812980181)
$restoredC = @(201.2072222222222,204.4562962962963,205.0889506172839,203.5841563786008,204.3764677640603,204.349858253315,204.1034941319921,204.2766067164558,204.2433197005876,204.2078068496785,204.242577755574)

for ($i = 0; $i -lt $restoredB.Length; $i++) {
    $row = $i + 52
    $ws.Cells.Item($row, 2).Value = $restoredB[$i]
    $ws.Cells.Item($row, 3).Value = $restoredC[$i]
}

Write-Host "edit applied"
